$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.235.34"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.905.80"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3957"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07965"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9999"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "1.936.03"
$ws.Range("E13").Value = "  +4.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.130"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.769"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06967"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "29.252.07"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.356"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "2.136.04"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.046"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.909"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09391"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9233"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.350"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.349"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.262"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05831"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.171"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5756"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.991"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5420"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.220"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.878"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.580"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("E51").Value = "  -5.41%  "
